$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mon Dec 11 14_56_07 2023"
$ws.Range("B2").Value = "loc"
$ws.Range("C2").Value = 5

$ws.Range("A3").Value = "Mon Dec 11 14_58_36 2023"
$ws.Range("B3").Value = "loc"
$ws.Range("C3").Value = 10

$ws.Range("A4").Value = "Mon Dec 11 15_02_13 2023"
$ws.Range("B4").Value = "loc"
$ws.Range("C4").Value = 5
